$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 141949969.33
$ws.Range("P2").Value = 1372289.09
$ws.Range("Q2").Value = 70656519.59999999

$ws.Range("R2").Value = "'"
$ws.Range("R2").ClearFormats()

$ws.Range("S2").Value = 35232251.34

$ws.Range("T2").Value = "'"
$ws.Range("T2").ClearFormats()

$ws.Range("U2").Value = 22857516.57

$ws.Range("V2").Value = "'"
$ws.Range("V2").ClearFormats()

$ws.Range("W2").Value = 32848460.38
$ws.Range("X2").Value = 17502019.36

$ws.Range("Y2").Value = "'"
$ws.Range("Y2").ClearFormats()

$ws.Range("Z2").Value = 4015091.91

$ws.Range("AA2").Value = "'"
$ws.Range("AA2").ClearFormats()

$ws.Range("AB2").Value = 109101508.95

$ws.Range("AC2").Value = "'"
$ws.Range("AC2").ClearFormats()

$ws.Range("AD2").Value = "'"
$ws.Range("AD2").ClearFormats()

$ws.Range("AE2").Value = "'"
$ws.Range("AE2").ClearFormats()

$ws.Range("AF2").Value = 422.574935337
$ws.Range("AG2").Value = 23.1408717699
